$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" (G2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-15 23:01:13"

# zh-cn sheet: "Correspond Handoff Datetime" (H2) and "Correspond Handback DateTime" (K2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-15 23:01:08"
$wsZhCn.Range("K2").Value = "2016-08-15 23:01:34"

# de-de sheet: "Correspond Handoff Datetime" (H2) and "Correspond Handback DateTime" (K2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-15 23:01:13"
$wsDeDe.Range("K2").Value = "2016-08-15 23:01:42"
